# Apply the edit described by the diff to the "Directorio" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Directorio")

# Row 26: the node that used to duplicate the ID "5.1.1." actually
# represents "5.1.1.1." at level 3 (was mistakenly tagged level 4).
$ws.Range("A26").Value = "5.1.1.1."
$ws.Range("C26").Value = 3

# Row 27: new child node "5.1.1.2." whose parent ID is now a literal
# value equal to the (corrected) row 26 ID, instead of the old formula
# "=+A26".
$ws.Range("A27").Value = "5.1.1.2."
$ws.Range("B27").Value = "5.1.1.1."

# Update the active selection to D8, matching the saved view state.
$ws.Range("D8").Select()
